$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells for the season-record columns, matching
# the look of the existing header row (bold, bordered, centered)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-46) gets the same season record: 91 wins, 71 losses, 0 ties
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 91
    $ws.Cells.Item($row, 31).Value = 71
    $ws.Cells.Item($row, 32).Value = 0
}
